$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1) cells:
#    "<name>_old" -> "<name>_FV2310"
#    "<name>_new" -> "<name>_FV2404"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value2 = ($val -replace "_old$", "_FV2310")
        } elseif ($val -like "*_new") {
            $cell.Value2 = ($val -replace "_new$", "_FV2404")
        }
    }
}

# 2. Freeze the header row (freeze pane below row 1).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the used range into an Excel Table ("Table1") with a header row,
#    using the already-renamed header names as the table's column headers.
$rng = $ws.Range("A1:U81")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
